$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (Coin name / Link) - plain text, no numeric coercion risk
$textUpdates = @{
    'B6' = 'FTXToken'
    'C6' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'B7' = 'MXToken'
    'C7' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'B8' = 'LiechtensteinCryptoassetsExchange'
    'C8' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'B9' = 'WazirX'
    'C9' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'B10' = 'MandalaExchangeToken'
    'C10' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'B11' = 'BitrueCoin'
    'C11' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'B12' = 'BitMartToken'
    'C12' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'B13' = 'BitForexToken'
    'C13' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'B14' = 'TigerCash'
    'C14' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'B15' = 'UpBots'
    'C15' = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'B17' = 'GateToken'
    'C17' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Numeric-looking columns (Price / Volume%) - must stay as literal text,
# matching original inline-string cells, so force text format before
# assignment and reset the style back to Normal afterwards so no visual
# / style artifacts remain on the cell.
$textForcedUpdates = @{
    'D2' = '303.12'
    'E2' = '2.09%'
    'D3' = '44.15'
    'E3' = '7.08%'
    'D4' = '5.106'
    'E4' = '1.96%'
    'D5' = '0.07739'
    'E5' = '3.21%'
    'D6' = '1.617'
    'E6' = '2.85%'
    'D7' = '1.046'
    'E7' = '12.95%'
    'D8' = '0.1284'
    'E8' = '5.15%'
    'D9' = '0.1863'
    'E9' = '1.36%'
    'D10' = '0.09235'
    'E10' = '4.51%'
    'D11' = '0.04152'
    'E11' = '1.53%'
    'D12' = '0.1047'
    'E12' = '-0.67%'
    'D13' = '0.001281'
    'E13' = '-0.30%'
    'D14' = '0.005745'
    'E14' = '-1.33%'
    'D15' = '0.007489'
    'E15' = '1,911.15%'
    'D16' = '3.355'
    'E16' = '0.29%'
    'D17' = '4.413'
    'E17' = '1.14%'
    'D19' = '0.3354'
    'E19' = '2.01%'
    'D20' = '8.037'
    'E20' = '0.30%'
    'E21' = '-2.97%'
    'E22' = '7.22%'
    'D23' = '0.04196'
    'E23' = '3.82%'
    'E24' = '0.97%'
    'D25' = '0.004406'
    'E25' = '13.47%'
    'E26' = '9.49%'
    'D38' = '0.02508'
    'E38' = '3.90%'
    'D39' = '0.05306'
    'E39' = '1.62%'
    'D40' = '0.005813'
    'E40' = '-2.97%'
    'D41' = '0.007716'
    'E41' = '-0.98%'
    'D42' = '0.1355'
    'E42' = '2.26%'
    'D43' = '0.007353'
    'E43' = '-0.17%'
    'D44' = '0.007511'
    'E44' = '-7.49%'
    'D45' = '0.3017'
    'E45' = '1.63%'
    'D46' = '0.00006674'
    'E46' = '6.36%'
    'E47' = '-0.24%'
    'D48' = '0.04344'
    'E48' = '-3.69%'
    'E49' = '-0.24%'
    'D50' = '0.0001995'
    'E50' = '-0.24%'
}
foreach ($addr in $textForcedUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$addr]
    $cell.Style = "Normal"
}
